$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 14: fill in the values (order chosen to match shared-string allocation order) ---
$ws.Range("C14").Value = "15/12/2018"
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = "PMO"
$ws.Range("G14").Value = "Creación G.Costes y Justificación"
$ws.Range("H14").Value = "Todos"
$ws.Range("I14").Value = "Creación G.Costes y Justificación"
$ws.Range("J14").Value = "-"
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = "PMO"
$ws.Range("M14").Value = "PMO"

# --- Dates for rows 15 and 16 ---
$ws.Range("C15").Value = "23/12/2018"
$ws.Range("C16").Value = "24/12/2018"

# --- Descriptions, row 16 before row 15 (matches original authoring order) ---
$ws.Range("G16").Value = "Correción documento costes"
$ws.Range("I16").Value = "Correción documento costes"

$ws.Range("G15").Value = "Adición jutificación"
$ws.Range("I15").Value = "Adición jutificación"

# --- Remaining row 15 values ---
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1.1
$ws.Range("F15").Value = "PMO"
$ws.Range("H15").Value = "Todos"
$ws.Range("J15").Value = "-"
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = "PMO"
$ws.Range("M15").Value = "PMO"

# --- Remaining row 16 values ---
$ws.Range("D16").Value = 1.1
$ws.Range("E16").Value = 1.2
$ws.Range("F16").Value = "Dpto.ctrl y Calidad"
$ws.Range("H16").Value = "Todos"
$ws.Range("J16").Value = "-"
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = "Dpto.ctrl y Calidad"
$ws.Range("M16").Value = "Dpto.ctrl y Calidad"

# --- Alignment fix-ups for rows 15 & 16 (center horizontally, like row 13) ---
$ws.Range("D15:F15").HorizontalAlignment = -4108
$ws.Range("J15:M15").HorizontalAlignment = -4108
$ws.Range("D16:F16").HorizontalAlignment = -4108
$ws.Range("J16:M16").HorizontalAlignment = -4108

# I15/I16 revert from left-aligned to general alignment
$ws.Range("I15").HorizontalAlignment = 1
$ws.Range("I16").HorizontalAlignment = 1

# --- Selection as left by the editing session ---
$ws.Range("J15:M16").Select()
